# Updates the cryptos price/volume table (columns D and E, rows 2-51).
# Note: several "Price" values in column D look numeric (e.g. "526.57") but
# must stay plain text, matching the original workbook. Assigning such a
# string directly to Range.Value would make Excel auto-convert it to a
# number, so for those cells we instead write a text-forced value
# (leading apostrophe) into a scratch cell (Z1) and copy/paste-special
# just the value into the target cell, then clean up the scratch cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.774.87"
$ws.Range("E2").Value = "  -4.30%  "
$ws.Range("D3").Value = "3.163.03"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("Z1").Value = "'526.57"
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -6.03%  "
$ws.Range("Z1").Value = "'133.52"
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -7.39%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.162.42"
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("E9").Value = "  -6.25%  "
$ws.Range("Z1").Value = "'7.35"
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("Z1").Value = "'0.111"
$ws.Range("Z1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -7.31%  "
$ws.Range("Z1").Value = "'0.391"
$ws.Range("Z1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "3.708.80"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("Z1").Value = "'25.87"
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -5.15%  "
$ws.Range("D16").Value = "3.170.86"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").Value = "57.809.21"
$ws.Range("E17").Value = "  -4.22%  "
$ws.Range("Z1").Value = "'0.0000153"
$ws.Range("Z1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -7.99%  "
$ws.Range("Z1").Value = "'5.83"
$ws.Range("Z1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("Z1").Value = "'13.09"
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -9.09%  "
$ws.Range("Z1").Value = "'8.07"
$ws.Range("Z1").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("Z1").Value = "'346.53"
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -7.22%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("Z1").Value = "'69.68"
$ws.Range("Z1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -6.02%  "
$ws.Range("Z1").Value = "'0.513"
$ws.Range("Z1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -6.86%  "
$ws.Range("D26").Value = "3.297.61"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").Value = "0.0₃0961"
$ws.Range("E27").Value = "  -8.83%  "
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("Z1").Value = "'1.00"
$ws.Range("Z1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("Z1").Value = "'6.89"
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("Z1").Value = "'0.997"
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("Z1").Value = "'1.88"
$ws.Range("Z1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -8.09%  "
$ws.Range("Z1").Value = "'6.92"
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -9.37%  "
$ws.Range("Z1").Value = "'21.76"
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("Z1").Value = "'4.91"
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -5.20%  "
$ws.Range("Z1").Value = "'160.16"
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("Z1").Value = "'6.26"
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("E39").Value = "  -7.85%  "
$ws.Range("Z1").Value = "'25.96"
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("Z1").Value = "'0.0698"
$ws.Range("Z1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("D42").Value = "3.191.81"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("Z1").Value = "'40.62"
$ws.Range("Z1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  -6.97%  "
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("Z1").Value = "'3.96"
$ws.Range("Z1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -6.06%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -8.27%  "
$ws.Range("D49").Value = "2.274.04"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("Z1").Value = "'6.19"
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -6.14%  "
$ws.Range("Z1").Value = "'20.55"
$ws.Range("Z1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -5.27%  "
$ws.Range("Z1").Clear()
